# 自动更新Excel文件 - 2025-12-26 23:13:31
# For every data row, recompute the "剩余" (remaining) column E based on the
# "总天" (total days) column D and the "开始时间" (start date, yyyymmdd) column F,
# counting down to a reference date of 2025-12-27.
# If the remaining count would drop to zero or below, the cycle is treated as
# renewed: the start date F is reset to the reference date (2025-12-27) and the
# remaining count E is reset back to the full total D.

function Get-JDN($y, $m, $d) {
    $a = [math]::Floor((14 - $m) / 12)
    $y2 = $y + 4800 - $a
    $m2 = $m + 12 * $a - 3
    $jdn = $d + [math]::Floor((153 * $m2 + 2) / 5) + 365 * $y2 + [math]::Floor($y2 / 4) - [math]::Floor($y2 / 100) + [math]::Floor($y2 / 400) - 32045
    return $jdn
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$todayY = 2025
$todayM = 12
$todayD = 27
$todayJDN = Get-JDN $todayY $todayM $todayD
$todaySerial = 20251227

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $d = $ws.Cells.Item($r, 4).Value2
    $f = $ws.Cells.Item($r, 6).Value2

    if ($d -eq $null) { continue }
    if ($f -eq $null) { continue }

    $fStr = [string]$f
    if ($fStr.Length -ne 8) { continue }

    $fy = [math]::Floor($f / 10000)
    $frem = $f % 10000
    $fm = [math]::Floor($frem / 100)
    $fd = $frem % 100

    $fJDN = Get-JDN $fy $fm $fd
    $elapsed = $todayJDN - $fJDN
    $newE = $d - $elapsed

    if ($newE -le 0) {
        $newE = $d
        $ws.Cells.Item($r, 6).Value2 = $todaySerial
    }

    $ws.Cells.Item($r, 5).Value2 = $newE
}
